# Exigence-fonctionnel.xlsx edit: assign team members (Répartition) to the
# first few EF_001/EF_002 requirement rows, tweak column B width and move
# the active selection to C9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Répartition") assignments for rows 2-9
$ws.Range("C2").Value = "Alex"
$ws.Range("C3").Value = "Alex"
$ws.Range("C4").Value = "Alex"
$ws.Range("C5").Value = "Thomas"
$ws.Range("C6").Value = "Thomas"
$ws.Range("C7").Value = "Thomas"
$ws.Range("C8").Value = "Marion"
$ws.Range("C9").Value = "Marion"

# Column B got slightly narrower (manual resize by the author)
$ws.Columns("B").ColumnWidth = 128.66666666666669

# Active cell ends up on C9 after the edits
$ws.Range("C9").Select()
